# The "payment_frequency" column (G) is no longer part of the students
# import template - the header and the sample values ("annuel") are
# cleared out, while the column itself (its width/position) is left in
# place. This also naturally drops the now-unused "payment_frequency"
# and "annuel" entries from the shared-strings table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1:G8").ClearContents()

# Restore the cursor/selection to its position in the refreshed template.
[void]$ws.Range("G10").Select()
